$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.027.62"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "2.402.60"
$ws.Range("E3").Value = "  +1.27%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "507.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.02%  "
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("E8").Value = "  -0.08%  "
$ws.Range("D9").Value = "2.411.67"
$ws.Range("E9").Value = "  +0.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0978"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.71%  "
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("E13").Value = "  -3.05%  "
$ws.Range("D14").Value = "2.829.02"
$ws.Range("E14").Value = "  +1.29%  "
$ws.Range("D15").Value = "56.959.45"
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.81"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.83%  "
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("D18").Value = "2.386.03"
$ws.Range("E18").Value = "  -0.14%  "
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("E21").Value = "  -0.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.26"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.95%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.74%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.70"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +4.08%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.375"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.18%  "
$ws.Range("E28").Value = "  -0.23%  "
$ws.Range("E29").Value = "  +2.56%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "175.94"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.69%  "
$ws.Range("D31").Value = "0.0₃0727"
$ws.Range("E31").Value = "  +1.62%  "
$ws.Range("E32").Value = "  -0.38%  "
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.88"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.72%  "
$ws.Range("E35").Value = "  +0.17%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.996"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.24%  "
$ws.Range("E37").Value = "  +1.25%  "
$ws.Range("E38").Value = "  -0.19%  "
$ws.Range("E39").Value = "  +3.15%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.835"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.19%  "
$ws.Range("E42").Value = "  +0.57%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "132.73"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.48%  "
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.571"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("E47").Value = "  +1.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "251.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0490"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.31%  "
$ws.Range("E50").Value = "  +2.06%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.08"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +7.27%  "
